# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
# timestamps and sets the "Priority" column to "ht" for the batch of
# files handed off together (rows 9-14 on each report sheet).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$rows = 9..14

foreach ($r in $rows) {
    # Overview sheet: column G = Latest HO Xliff Generate Date.
    # This text is the same shared string as de-de's column H below,
    # so both must be updated together.
    $overview.Range("G$r").Value = "2016-08-31 10:23:32"

    # zh-cn sheet: column H = Latest Handoff Datetime, column E = Priority
    $zhcn.Range("H$r").Value = "2016-08-31 10:23:28"
    $zhcn.Range("E$r").Value = "ht"

    # de-de sheet: column H shares the same underlying text as Overview's
    # column G (same shared-string entry), column E = Priority
    $dede.Range("H$r").Value = "2016-08-31 10:23:32"
    $dede.Range("E$r").Value = "ht"
}
